# Auto-generated script applying odds updates from diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("G4").Value = 1.9
$ws.Range("I4").Value = 4.75
$ws.Range("J4").Value = 2.75
$ws.Range("M4").Value = 1.13
$ws.Range("N4").Value = 6
$ws.Range("X4").Value = 7.5
$ws.Range("Y4").Value = 10
$ws.Range("AE4").Value = 21
$ws.Range("AF4").Value = 81
$ws.Range("AG4").Value = 9
$ws.Range("AH4").Value = 21
$ws.Range("AI4").Value = 17
$ws.Range("AK4").Value = 41
$ws.Range("AW4").Value = 6
$ws.Range("AZ4").Value = 101
$ws.Range("BA4").Value = 151

# Row 6
$ws.Range("M6").Value = 1.04
$ws.Range("N6").Value = 13
$ws.Range("Q6").Value = 1.83
$ws.Range("R6").Value = 2.03

# Row 7
$ws.Range("H7").Value = 2.85
$ws.Range("I7").Value = 3.55
$ws.Range("L7").Value = 4
$ws.Range("T7").Value = 2.4
$ws.Range("AA7").Value = 20
$ws.Range("AD7").Value = 5.6
$ws.Range("AG7").Value = 9.25
$ws.Range("AH7").Value = 19
$ws.Range("AL7").Value = 45
$ws.Range("AN7").Value = 3.95
$ws.Range("AT7").Value = 2.37
$ws.Range("AX7").Value = 20

# Row 9
$ws.Range("M9").Value = 1.08
$ws.Range("O9").Value = 1.4

# Row 10
$ws.Range("G10").Value = 2.5
$ws.Range("I10").Value = 3.4
$ws.Range("J10").Value = 3.4
$ws.Range("K10").Value = 1.83
$ws.Range("L10").Value = 4
$ws.Range("M10").Value = 1.14
$ws.Range("N10").Value = 5.5
$ws.Range("AA10").Value = 26
$ws.Range("AI10").Value = 13
$ws.Range("AJ10").Value = 34
$ws.Range("AP10").Value = 34

# Row 11
$ws.Range("G11").Value = 2.05
$ws.Range("I11").Value = 3.9
$ws.Range("K11").Value = 2.1
$ws.Range("L11").Value = 4
$ws.Range("U11").Value = 1.8
$ws.Range("V11").Value = 1.91
$ws.Range("W11").Value = 7.5
$ws.Range("X11").Value = 9.5
$ws.Range("Z11").Value = 17
$ws.Range("AC11").Value = 9.5
$ws.Range("AF11").Value = 51
$ws.Range("AH11").Value = 19
$ws.Range("AI11").Value = 13
$ws.Range("AJ11").Value = 41
$ws.Range("AK11").Value = 29
$ws.Range("AL11").Value = 41
$ws.Range("AM11").Value = 251
$ws.Range("AN11").Value = 4
$ws.Range("AO11").Value = 11
$ws.Range("AX11").Value = 21
$ws.Range("AY11").Value = 29
$ws.Range("AZ11").Value = 67
$ws.Range("BB11").Value = 201

# Row 12
$ws.Range("G12").Value = 1.25

# Row 14
$ws.Range("H14").Value = 2.92
$ws.Range("J14").Value = 3.5
$ws.Range("L14").Value = 3
$ws.Range("M14").Value = 1.02
$ws.Range("N14").Value = 7.1
$ws.Range("O14").Value = 1.36
$ws.Range("P14").Value = 2.67
$ws.Range("Q14").Value = 2.05
$ws.Range("U14").Value = 1.75
$ws.Range("V14").Value = 1.85
$ws.Range("AA14").Value = 28
$ws.Range("AB14").Value = 37
$ws.Range("AC14").Value = 7.8
$ws.Range("AG14").Value = 7.2
$ws.Range("AH14").Value = 11.5
$ws.Range("AI14").Value = 9.25
$ws.Range("AK14").Value = 21
$ws.Range("AL14").Value = 32
$ws.Range("AN14").Value = 4.85
$ws.Range("AO14").Value = 16.5
$ws.Range("AP14").Value = 23
$ws.Range("AQ14").Value = 80
$ws.Range("AR14").Value = 110
$ws.Range("AS14").Value = 300
$ws.Range("AT14").Value = 2.45
$ws.Range("AU14").Value = 6.6
$ws.Range("AV14").Value = 55
$ws.Range("AY14").Value = 20
$ws.Range("BA14").Value = 90
$ws.Range("BB14").Value = 250
